$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to S.N. 8, Date 45340 (Feb 18, 2024); C9 was empty.
$ws.Range("C9").Value = "Text formatting and tables"
